$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the two header cells (K1, L1). This also renames the underlying
# Table1 ListColumns ("authentic_source_id" -> "starting_date",
# "authentic_source_name" -> "ending_date") since the table's headers live
# in these cells.
$ws.Range("K1").Value = "starting_date"
$ws.Range("L1").Value = "ending_date"

# Replace the authentic_source_id / authentic_source_name values in both
# data rows with starting_date / ending_date date values (real date
# serials, not text), formatted with a date number format.
# Format K2 first, then copy its formatting onto L2, K3 and L3 so all four
# cells consistently share the resulting date style.
$ws.Range("K2").Value2 = 43863
$ws.Range("K2").NumberFormat = "mm-dd-yy"
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("L3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("L2").Value2 = 45690
$ws.Range("K3").Value2 = 43863
$ws.Range("L3").Value2 = 45690

# Update the view: scroll so column H is the left-most visible column and
# move the active selection to L4.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("L4").Select()
